$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.109.27'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '3.614.67'
$ws.Range('E3').Value = '  +3.27%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''604.52'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').Value = '''195.30'
$ws.Range('E6').Value = '  -1.32%  '
$ws.Range('D7').Value = '''0.627'
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -2.09%  '
$ws.Range('D10').Value = '''0.653'
$ws.Range('E10').Value = '  -0.11%  '
$ws.Range('D11').Value = '''54.06'
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('E12').Value = '  +0.70%  '
$ws.Range('D13').Value = '''9.55'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').Value = '4.193.63'
$ws.Range('E14').Value = '  +3.34%  '
$ws.Range('D15').Value = '''13.23'
$ws.Range('E15').Value = '  +4.83%  '
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = '''19.23'
$ws.Range('E17').Value = '  +0.86%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '70.294.70'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.614.16'
$ws.Range('E19').Value = '  +3.24%  '
$ws.Range('E20').Value = '  +1.65%  '
$ws.Range('D21').Value = '''0.996'
$ws.Range('E21').Value = '  +0.60%  '
$ws.Range('D22').Value = '''17.82'
$ws.Range('E22').Value = '  -2.82%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Value = '''103.19'
$ws.Range('E23').Value = '  -1.35%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = '''5.16'
$ws.Range('E24').Value = '  +2.61%  '
$ws.Range('D25').Value = '''4.64'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('D26').Value = '''3.07'
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('D27').Value = '''10.80'
$ws.Range('E27').Value = '  -1.67%  '
$ws.Range('E28').Value = '  -1.47%  '
$ws.Range('D29').Value = '''34.02'
$ws.Range('E29').Value = '  +1.16%  '
$ws.Range('D30').Value = '''4.44'
$ws.Range('E30').Value = '  -1.81%  '
$ws.Range('D31').Value = '''7.17'
$ws.Range('E31').Value = '  -1.16%  '
$ws.Range('D32').Value = '''12.37'
$ws.Range('E32').Value = '  -2.85%  '
$ws.Range('E33').Value = '  +1.18%  '
$ws.Range('D34').Value = '''63.25'
$ws.Range('D35').Value = '0.0₃0874'
$ws.Range('E35').Value = '  +8.04%  '
$ws.Range('D36').Value = '3.938.07'
$ws.Range('E36').Value = '  +5.52%  '
$ws.Range('E37').Value = '  +7.19%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '''523.19'
$ws.Range('E38').Value = '  +2.51%  '
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').Value = '''1.00'
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('D40').Value = '''37.23'
$ws.Range('E40').Value = '  +1.40%  '
$ws.Range('E41').Value = '  +0.70%  '
$ws.Range('D42').Value = '''3.56'
$ws.Range('E42').Value = '  +1.07%  '
$ws.Range('E43').Value = '  -2.17%  '
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').Value = '''2.86'
$ws.Range('E45').Value = '  +0.87%  '
$ws.Range('E46').Value = '  +1.01%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '''3.34'
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('E48').Value = '  -1.45%  '
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('E50').Value = '  +4.09%  '
$ws.Range('E51').Value = '  +3.42%  '
